# Update Daily Report: 2026-02-04
# Adds a new day's data (date serial 46056) to Daily_Data, and rolls the
# resulting totals up into Today_Summary and Monthly_Stats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append 22 new rows (rows 464-485) for date 46056
# ---------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @('ASAHI DEPOSITORY LLC Registered', 0, 0, 0, 0, 0, 0),
    @('ASAHI DEPOSITORY LLC Eligible', 0, 0, 0, 0, 0, 0),
    @('BRINK''S, INC. Registered', 76497.842, 0, 0, 0, 0, 76497.842),
    @('BRINK''S, INC. Eligible', 42030.257, 56826.488, 0, 56826.488, 0, 98856.745),
    @('CNT DEPOSITORY, INC. Registered', 1246.06, 0, 0, 0, 0, 1246.06),
    @('CNT DEPOSITORY, INC. Eligible', 0, 0, 0, 0, 0, 0),
    @('DELAWARE DEPOSITORY Registered', 1633.941, 0, 0, 0, 0, 1633.941),
    @('DELAWARE DEPOSITORY Eligible', 18459.584, 0, 0, 0, 0, 18459.584),
    @('HSBC BANK, USA Registered', 1394.758, 0, 0, 0, 0, 1394.758),
    @('HSBC BANK, USA Eligible', 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 2395.448, 0, 0, 0, 0, 2395.448),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 0, 0, 0, 0, 0, 0),
    @('JP MORGAN CHASE BANK NA Registered', 114985.579, 0, 0, 0, 0, 114985.579),
    @('JP MORGAN CHASE BANK NA Eligible', 75484.511, 0, 0, 0, 0, 75484.511),
    @('LOOMIS INTERNATIONAL (US) LLC Registered', 63745.991, 0, 0, 0, 0, 63745.991),
    @('LOOMIS INTERNATIONAL (US) LLC Eligible', 132077.206, 0, 0, 0, 0, 132077.206),
    @('MALCA-AMIT USA, LLC Registered', 395.145, 0, 0, 0, 0, 395.145),
    @('MALCA-AMIT USA, LLC Eligible', 0, 0, 0, 0, 0, 0),
    @('MANFRA, TORDELLA & BROOKES, LLC Registered', 50220.42, 0, 0, 0, 0, 50220.42),
    @('MANFRA, TORDELLA & BROOKES, LLC Eligible', 1271.373, 0, 0, 0, 0, 1271.373),
    @('STONEX PRECIOUS METALS LLC Registered', 14122.765, 0, 0, 0, 0, 14122.765),
    @('STONEX PRECIOUS METALS LLC Eligible', 16.075, 0, 0, 0, 0, 16.075)
)

$startRow = 464
$newDate = 46056

$r = $startRow
foreach ($row in $newRows) {
    $wsDaily.Cells.Item($r, 1).Value = $newDate
    $wsDaily.Cells.Item($r, 1).NumberFormat = $wsDaily.Cells.Item($r - 1, 1).NumberFormat
    $wsDaily.Cells.Item($r, 2).Value = $row[0]
    $wsDaily.Cells.Item($r, 3).Value = $row[1]
    $wsDaily.Cells.Item($r, 4).Value = $row[2]
    $wsDaily.Cells.Item($r, 5).Value = $row[3]
    $wsDaily.Cells.Item($r, 6).Value = $row[4]
    $wsDaily.Cells.Item($r, 7).Value = $row[5]
    $wsDaily.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Today_Summary: BRINK'S, INC. row (row 3) - Eligible & Total_Stock
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("Today_Summary")
$wsToday.Range("B3").Value = 98856.745
$wsToday.Range("D3").Value = 175354.587

# ---------------------------------------------------------------------
# 3) Monthly_Stats: 2026-02 month totals (row 2) - Eligible & Grand_Total
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")
$wsMonthly.Range("B2").Value = 335447.473
$wsMonthly.Range("D2").Value = 662085.422

# ---------------------------------------------------------------------
# 4) Monthly_Stats: BRINK'S, INC. Eligible detail row (row 10) for 2026-02
#    - RECEIVED & TOTAL_TODAY
# ---------------------------------------------------------------------
$wsMonthly.Range("C10").Value = 56826.488
$wsMonthly.Range("E10").Value = 98856.745

Write-Output "Daily report updated for 2026-02-04 (46056)"
